$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "ManufacturerId"
$ws.Range("C1").Value = "SimCountryCode"
$ws.Range("D1").Value = "SimNumber"
$ws.Range("E1").Value = "OwnerPartyType"
$ws.Range("F1").Value = "OwnerId"
$ws.Range("G1").Value = "AssignedPartyType"
$ws.Range("H1").Value = "AssignedId"

# Columns whose header text is longer than the default width get
# widened (as Excel does on "best fit" auto-size of a column).
$ws.Columns.Item(2).ColumnWidth = 12.9167
$ws.Columns.Item(3).ColumnWidth = 13.5834
$ws.Columns.Item(5).ColumnWidth = 13.7501
$ws.Columns.Item(7).ColumnWidth = 15.7501

$ws.Range("H12").Select()
